# Reshuffle the per-observation data in rows 7-21 of the "Artfynd" sheet.
# The underlying records (species/coordinates/etc.) were re-sorted; each
# destination row ends up showing the data that used to live in a
# different source row, per the mapping below (row 13 is unchanged).
#
# Columns that actually carry per-record data: A,B,D,E,F,G,H,Q,R plus the
# sparse K,L,M,N (placeholder, always blank) and AC (free-text comment,
# only present on "Tretåig hackspett" rows in this sheet).
# Columns C,I,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY are identical
# across all of these rows, so they do not need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 7
$lastRow = 21

# destination row -> source row (where its new content currently lives)
$mapping = @{
    7  = 15
    8  = 18
    9  = 14
    10 = 9
    11 = 10
    12 = 11
    13 = 13
    14 = 12
    15 = 8
    16 = 20
    17 = 21
    18 = 7
    19 = 16
    20 = 19
    21 = 17
}

# Column letters -> 1-based column index
$colA  = 1    # Id
$colB  = 2    # Taxonsorteringsordning
$colD  = 4    # Rödlistade
$colE  = 5    # TaxonId
$colF  = 6    # Artnamn
$colG  = 7    # Vetenskapligt namn
$colH  = 8    # Auktor
$colK  = 11   # Ålder-Stadium (sparse placeholder)
$colL  = 12   # Kön (sparse placeholder)
$colM  = 13   # Aktivitet (sparse placeholder)
$colN  = 14   # Metod (sparse placeholder)
$colQ  = 17   # Ost
$colR  = 18   # Nord
$colAC = 29   # Publik kommentar

# 1) Snapshot every source row's current values BEFORE writing anything,
#    so the permutation is computed from the original state.
$saved = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $saved[$r] = @{
        A  = $ws.Cells.Item($r, $colA).Value()
        B  = $ws.Cells.Item($r, $colB).Value()
        D  = $ws.Cells.Item($r, $colD).Value()
        E  = $ws.Cells.Item($r, $colE).Value()
        F  = $ws.Cells.Item($r, $colF).Value()
        G  = $ws.Cells.Item($r, $colG).Value()
        H  = $ws.Cells.Item($r, $colH).Value()
        K  = $ws.Cells.Item($r, $colK).Value()
        L  = $ws.Cells.Item($r, $colL).Value()
        M  = $ws.Cells.Item($r, $colM).Value()
        N  = $ws.Cells.Item($r, $colN).Value()
        Q  = $ws.Cells.Item($r, $colQ).Value()
        R  = $ws.Cells.Item($r, $colR).Value()
        AC = $ws.Cells.Item($r, $colAC).Value()
    }
}

# 2) Write each destination row's cells from its mapped source row's
#    snapshot.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $mapping[$r]
    $data = $saved[$src]

    $ws.Cells.Item($r, $colA).Value = $data.A
    $ws.Cells.Item($r, $colB).Value = $data.B
    $ws.Cells.Item($r, $colD).Value = $data.D
    $ws.Cells.Item($r, $colE).Value = $data.E
    $ws.Cells.Item($r, $colF).Value = $data.F
    $ws.Cells.Item($r, $colG).Value = $data.G
    $ws.Cells.Item($r, $colH).Value = $data.H
    $ws.Cells.Item($r, $colK).Value = $data.K
    $ws.Cells.Item($r, $colL).Value = $data.L
    $ws.Cells.Item($r, $colM).Value = $data.M
    $ws.Cells.Item($r, $colN).Value = $data.N
    $ws.Cells.Item($r, $colQ).Value = $data.Q
    $ws.Cells.Item($r, $colR).Value = $data.R
    $ws.Cells.Item($r, $colAC).Value = $data.AC
}
